# Applies the GitHub Actions IESO report refresh captured in the
# 2025-06-08T15:07:43 run: refreshed CreatedAt banner plus updated
# Predispatch Hourly Energy LMP figures (columns S:Z) across the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-06-08T15:07:43"

# Row 4
$ws.Cells.Item(4, 19).Value = 156.03
$ws.Cells.Item(4, 20).Value = 138.27
$ws.Cells.Item(4, 21).Value = 33.98
$ws.Cells.Item(4, 22).Value = 50.9
$ws.Cells.Item(4, 23).Value = 36.19
$ws.Cells.Item(4, 24).Value = 36.22
$ws.Cells.Item(4, 25).Value = 36.04
$ws.Cells.Item(4, 26).Value = 14.79

# Row 6
$ws.Cells.Item(6, 19).Value = -2.34
$ws.Cells.Item(6, 20).Value = -2.07
$ws.Cells.Item(6, 22).Value = -0.51
$ws.Cells.Item(6, 23).Value = -0.43
$ws.Cells.Item(6, 24).Value = -0.25
$ws.Cells.Item(6, 25).Value = -0.18
$ws.Cells.Item(6, 26).Value = 0.03

# Row 9
$ws.Cells.Item(9, 19).Value = 149.27
$ws.Cells.Item(9, 20).Value = 137.18
$ws.Cells.Item(9, 21).Value = 34.11
$ws.Cells.Item(9, 22).Value = 51.72
$ws.Cells.Item(9, 23).Value = 36.33
$ws.Cells.Item(9, 24).Value = 36.81
$ws.Cells.Item(9, 25).Value = 37.18
$ws.Cells.Item(9, 26).Value = 15.51

# Row 11
$ws.Cells.Item(11, 19).Value = -9.109999999999999
$ws.Cells.Item(11, 20).Value = -3.16
$ws.Cells.Item(11, 23).Value = -0.29
$ws.Cells.Item(11, 24).Value = 0.33
$ws.Cells.Item(11, 25).Value = 0.97
$ws.Cells.Item(11, 26).Value = 0.74

# Row 14
$ws.Cells.Item(14, 19).Value = 150.42
$ws.Cells.Item(14, 20).Value = 137.18
$ws.Cells.Item(14, 21).Value = 34.08
$ws.Cells.Item(14, 22).Value = 51.72
$ws.Cells.Item(14, 23).Value = 36.33
$ws.Cells.Item(14, 24).Value = 36.81
$ws.Cells.Item(14, 25).Value = 37.18
$ws.Cells.Item(14, 26).Value = 15.52

# Row 15
$ws.Cells.Item(15, 19).Value = 1.29
$ws.Cells.Item(15, 20).Value = 0

# Row 16
$ws.Cells.Item(16, 19).Value = -9.25
$ws.Cells.Item(16, 20).Value = -3.16
$ws.Cells.Item(16, 21).Value = -0.31
$ws.Cells.Item(16, 23).Value = -0.29
$ws.Cells.Item(16, 24).Value = 0.33
$ws.Cells.Item(16, 25).Value = 0.97
$ws.Cells.Item(16, 26).Value = 0.76

# Row 19
$ws.Cells.Item(19, 19).Value = 156.65
$ws.Cells.Item(19, 20).Value = 138.95
$ws.Cells.Item(19, 21).Value = 34.18
$ws.Cells.Item(19, 22).Value = 51.26
$ws.Cells.Item(19, 23).Value = 36.44
$ws.Cells.Item(19, 24).Value = 36.47
$ws.Cells.Item(19, 25).Value = 36.36
$ws.Cells.Item(19, 26).Value = 14.94

# Row 21
$ws.Cells.Item(21, 19).Value = -1.72
$ws.Cells.Item(21, 20).Value = -1.39
$ws.Cells.Item(21, 21).Value = -0.21
$ws.Cells.Item(21, 22).Value = -0.15
$ws.Cells.Item(21, 23).Value = -0.18
$ws.Cells.Item(21, 24).Value = 0
$ws.Cells.Item(21, 25).Value = 0.15
$ws.Cells.Item(21, 26).Value = 0.18

# Row 24
$ws.Cells.Item(24, 19).Value = 156.65
$ws.Cells.Item(24, 20).Value = 138.95
$ws.Cells.Item(24, 21).Value = 34.18
$ws.Cells.Item(24, 22).Value = 51.26
$ws.Cells.Item(24, 23).Value = 36.44
$ws.Cells.Item(24, 24).Value = 36.47
$ws.Cells.Item(24, 25).Value = 36.36
$ws.Cells.Item(24, 26).Value = 14.94

# Row 26
$ws.Cells.Item(26, 19).Value = -1.72
$ws.Cells.Item(26, 20).Value = -1.39
$ws.Cells.Item(26, 21).Value = -0.21
$ws.Cells.Item(26, 22).Value = -0.15
$ws.Cells.Item(26, 23).Value = -0.18
$ws.Cells.Item(26, 24).Value = 0
$ws.Cells.Item(26, 25).Value = 0.15
$ws.Cells.Item(26, 26).Value = 0.18

# Row 29
$ws.Cells.Item(29, 19).Value = 157.9
$ws.Cells.Item(29, 20).Value = 139.92
$ws.Cells.Item(29, 21).Value = 34.45
$ws.Cells.Item(29, 22).Value = 51.72
$ws.Cells.Item(29, 23).Value = 36.77
$ws.Cells.Item(29, 24).Value = 36.77
$ws.Cells.Item(29, 25).Value = 36.73
$ws.Cells.Item(29, 26).Value = 15.11

# Row 31
$ws.Cells.Item(31, 19).Value = -0.47
$ws.Cells.Item(31, 20).Value = -0.42
$ws.Cells.Item(31, 22).Value = 0.31
$ws.Cells.Item(31, 23).Value = 0.15
$ws.Cells.Item(31, 24).Value = 0.29
$ws.Cells.Item(31, 25).Value = 0.51
$ws.Cells.Item(31, 26).Value = 0.35

# Row 34
$ws.Cells.Item(34, 19).Value = 150.14
$ws.Cells.Item(34, 20).Value = 136.78
$ws.Cells.Item(34, 21).Value = 34.11
$ws.Cells.Item(34, 22).Value = 51.82
$ws.Cells.Item(34, 23).Value = 36.19
$ws.Cells.Item(34, 24).Value = 37.07
$ws.Cells.Item(34, 25).Value = 37.61
$ws.Cells.Item(34, 26).Value = 15.69

# Row 35
$ws.Cells.Item(35, 19).Value = 1.29
$ws.Cells.Item(35, 20).Value = 0

# Row 36
$ws.Cells.Item(36, 19).Value = -9.529999999999999
$ws.Cells.Item(36, 20).Value = -3.56
$ws.Cells.Item(36, 22).Value = 0.41
$ws.Cells.Item(36, 23).Value = -0.43
$ws.Cells.Item(36, 24).Value = 0.59
$ws.Cells.Item(36, 25).Value = 1.39
$ws.Cells.Item(36, 26).Value = 0.93

# Row 39
$ws.Cells.Item(39, 19).Value = 156.03
$ws.Cells.Item(39, 20).Value = 138.27
$ws.Cells.Item(39, 21).Value = 33.98
$ws.Cells.Item(39, 22).Value = 50.9
$ws.Cells.Item(39, 23).Value = 36.19
$ws.Cells.Item(39, 24).Value = 36.22
$ws.Cells.Item(39, 25).Value = 36.04
$ws.Cells.Item(39, 26).Value = 14.79

# Row 41
$ws.Cells.Item(41, 19).Value = -2.34
$ws.Cells.Item(41, 20).Value = -2.07
$ws.Cells.Item(41, 22).Value = -0.51
$ws.Cells.Item(41, 23).Value = -0.43
$ws.Cells.Item(41, 24).Value = -0.25
$ws.Cells.Item(41, 25).Value = -0.18
$ws.Cells.Item(41, 26).Value = 0.03

# Row 44
$ws.Cells.Item(44, 19).Value = 157.27
$ws.Cells.Item(44, 20).Value = 139.64
$ws.Cells.Item(44, 21).Value = 34.21
$ws.Cells.Item(44, 22).Value = 51.26
$ws.Cells.Item(44, 23).Value = 36.51
$ws.Cells.Item(44, 24).Value = 36.51
$ws.Cells.Item(44, 25).Value = 36.29
$ws.Cells.Item(44, 26).Value = 14.85

# Row 46
$ws.Cells.Item(46, 19).Value = -1.1
$ws.Cells.Item(46, 20).Value = -0.7
$ws.Cells.Item(46, 23).Value = -0.11
$ws.Cells.Item(46, 24).Value = 0.04
$ws.Cells.Item(46, 25).Value = 0.07000000000000001

# Row 49
$ws.Cells.Item(49, 19).Value = 137.36
$ws.Cells.Item(49, 20).Value = 121.72
$ws.Cells.Item(49, 21).Value = 30.14
$ws.Cells.Item(49, 22).Value = 48.14
$ws.Cells.Item(49, 23).Value = 34.19
$ws.Cells.Item(49, 24).Value = 35.72
$ws.Cells.Item(49, 25).Value = 35.68
$ws.Cells.Item(49, 26).Value = 14.66

# Row 51
$ws.Cells.Item(51, 19).Value = -21.02
$ws.Cells.Item(51, 20).Value = -18.62
$ws.Cells.Item(51, 21).Value = -4.25
$ws.Cells.Item(51, 23).Value = -2.43
$ws.Cells.Item(51, 24).Value = -0.75
$ws.Cells.Item(51, 25).Value = -0.54

# Row 54
$ws.Cells.Item(54, 19).Value = 147.74
$ws.Cells.Item(54, 20).Value = 131.4
$ws.Cells.Item(54, 21).Value = 32.23
$ws.Cells.Item(54, 22).Value = 48.77
$ws.Cells.Item(54, 23).Value = 35.49
$ws.Cells.Item(54, 24).Value = 35.34
$ws.Cells.Item(54, 25).Value = 35.51
$ws.Cells.Item(54, 26).Value = 14.59

# Row 56
$ws.Cells.Item(56, 19).Value = -10.64
$ws.Cells.Item(56, 20).Value = -8.94
$ws.Cells.Item(56, 21).Value = -2.16
$ws.Cells.Item(56, 22).Value = -2.63
$ws.Cells.Item(56, 23).Value = -1.14
$ws.Cells.Item(56, 24).Value = -1.13
$ws.Cells.Item(56, 25).Value = -0.71

# Row 59
$ws.Cells.Item(59, 19).Value = 161.94
$ws.Cells.Item(59, 20).Value = 143.35
$ws.Cells.Item(59, 21).Value = 35.09
$ws.Cells.Item(59, 22).Value = 52.62
$ws.Cells.Item(59, 23).Value = 37.56
$ws.Cells.Item(59, 24).Value = 37.56
$ws.Cells.Item(59, 25).Value = 37.38
$ws.Cells.Item(59, 26).Value = 15.28

# Row 61
$ws.Cells.Item(61, 19).Value = 3.56
$ws.Cells.Item(61, 20).Value = 3.01
$ws.Cells.Item(61, 24).Value = 1.09
$ws.Cells.Item(61, 25).Value = 1.16
$ws.Cells.Item(61, 26).Value = 0.52

# Row 64
$ws.Cells.Item(64, 19).Value = 164.46
$ws.Cells.Item(64, 20).Value = 145.58
$ws.Cells.Item(64, 21).Value = 35.63
$ws.Cells.Item(64, 22).Value = 53.44
$ws.Cells.Item(64, 23).Value = 38.11
$ws.Cells.Item(64, 24).Value = 38.07
$ws.Cells.Item(64, 25).Value = 37.69
$ws.Cells.Item(64, 26).Value = 15.47

# Row 66
$ws.Cells.Item(66, 19).Value = 6.08
$ws.Cells.Item(66, 20).Value = 5.24
$ws.Cells.Item(66, 21).Value = 1.25
$ws.Cells.Item(66, 23).Value = 1.49
$ws.Cells.Item(66, 24).Value = 1.6
$ws.Cells.Item(66, 25).Value = 1.47
$ws.Cells.Item(66, 26).Value = 0.71

# Row 69
$ws.Cells.Item(69, 19).Value = 165.49
$ws.Cells.Item(69, 20).Value = 146.19
$ws.Cells.Item(69, 21).Value = 35.82
$ws.Cells.Item(69, 22).Value = 53.77
$ws.Cells.Item(69, 23).Value = 38.47
$ws.Cells.Item(69, 24).Value = 38.78
$ws.Cells.Item(69, 25).Value = 38.37
$ws.Cells.Item(69, 26).Value = 15.69

# Row 71
$ws.Cells.Item(71, 19).Value = 7.12
$ws.Cells.Item(71, 20).Value = 5.85
$ws.Cells.Item(71, 21).Value = 1.43
$ws.Cells.Item(71, 22).Value = 2.37
$ws.Cells.Item(71, 23).Value = 1.85
$ws.Cells.Item(71, 24).Value = 2.08
$ws.Cells.Item(71, 25).Value = 2.15
$ws.Cells.Item(71, 26).Value = 0.93

# Row 72
$ws.Cells.Item(72, 24).Value = 0.22

# Row 74
$ws.Cells.Item(74, 19).Value = 158.37
$ws.Cells.Item(74, 20).Value = 140.34
$ws.Cells.Item(74, 21).Value = 34.38
$ws.Cells.Item(74, 22).Value = 51.41
$ws.Cells.Item(74, 23).Value = 36.62
$ws.Cells.Item(74, 24).Value = 36.47
$ws.Cells.Item(74, 25).Value = 36.22
$ws.Cells.Item(74, 26).Value = 14.76

# Row 79
$ws.Cells.Item(79, 19).Value = 158.37
$ws.Cells.Item(79, 20).Value = 140.34
$ws.Cells.Item(79, 21).Value = 34.38
$ws.Cells.Item(79, 22).Value = 51.41
$ws.Cells.Item(79, 23).Value = 36.62
$ws.Cells.Item(79, 24).Value = 36.47
$ws.Cells.Item(79, 25).Value = 36.22
$ws.Cells.Item(79, 26).Value = 14.76

# Row 84
$ws.Cells.Item(84, 19).Value = 142.94
$ws.Cells.Item(84, 20).Value = 127.47
$ws.Cells.Item(84, 21).Value = 31.46
$ws.Cells.Item(84, 22).Value = 47.51
$ws.Cells.Item(84, 23).Value = 35.49
$ws.Cells.Item(84, 24).Value = 35.21
$ws.Cells.Item(84, 25).Value = 35.4
$ws.Cells.Item(84, 26).Value = 14.53

# Row 86
$ws.Cells.Item(86, 19).Value = -15.44
$ws.Cells.Item(86, 20).Value = -12.87
$ws.Cells.Item(86, 21).Value = -2.93
$ws.Cells.Item(86, 22).Value = -3.9
$ws.Cells.Item(86, 23).Value = -1.14
$ws.Cells.Item(86, 24).Value = -1.27
$ws.Cells.Item(86, 25).Value = -0.8100000000000001

# Row 89
$ws.Cells.Item(89, 19).Value = 157.9
$ws.Cells.Item(89, 20).Value = 139.92
$ws.Cells.Item(89, 21).Value = 34.45
$ws.Cells.Item(89, 22).Value = 51.72
$ws.Cells.Item(89, 23).Value = 36.77
$ws.Cells.Item(89, 24).Value = 36.77
$ws.Cells.Item(89, 25).Value = 36.77
$ws.Cells.Item(89, 26).Value = 15.11

# Row 91
$ws.Cells.Item(91, 19).Value = -0.47
$ws.Cells.Item(91, 20).Value = -0.42
$ws.Cells.Item(91, 21).Value = 0.07000000000000001
$ws.Cells.Item(91, 22).Value = 0.31
$ws.Cells.Item(91, 23).Value = 0.15
$ws.Cells.Item(91, 24).Value = 0.29
$ws.Cells.Item(91, 25).Value = 0.55
$ws.Cells.Item(91, 26).Value = 0.35
